# Femacal de La Calera - Chirimoya: add a new week of prices (row block
# 198-200) ahead of the existing history, pushing every later row down by
# three positions (old 198-226 -> new 201-229).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new week: insert 3 blank rows at the top of the block.
$ws.Rows("198:200").Insert()

# Seed the new rows from the (now shifted-down) row that used to be 198-200,
# so all the formatting/styles and the columns that don't change
# (B, C, E, F, G, H, I, J, K, L, Q, T, ...) come along for free.
$ws.Range("A201:T203").Copy()
$ws.Range("A198").PasteSpecial()

# Row 198 - Especial, Provincia del Elquí, week of 2022-09-22
$ws.Range("D198").Value2 = 44826
$ws.Range("M198").Value2 = 56
$ws.Range("N198").Value2 = 32000
$ws.Range("O198").Value2 = 32000
$ws.Range("P198").Value2 = 32000
$ws.Range("R198").Value2 = "Provincia del Elquí"
$ws.Range("S198").Value2 = 3200

# Row 199 - Primera, Provincia del Elquí, week of 2022-09-22
$ws.Range("D199").Value2 = 44826
$ws.Range("M199").Value2 = 67
$ws.Range("N199").Value2 = 30000
$ws.Range("O199").Value2 = 30000
$ws.Range("P199").Value2 = 30000
$ws.Range("R199").Value2 = "Provincia del Elquí"
$ws.Range("S199").Value2 = 3000

# Row 200 - Segunda, Provincia del Elquí, week of 2022-09-22
$ws.Range("D200").Value2 = 44826
$ws.Range("M200").Value2 = 60
$ws.Range("N200").Value2 = 27000
$ws.Range("O200").Value2 = 27000
$ws.Range("P200").Value2 = 27000
$ws.Range("R200").Value2 = "Provincia del Elquí"
$ws.Range("S200").Value2 = 2700
